$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna2"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.408030333333333
$ws.Range("H2").Value = 4.224091
$ws.Range("I2").Value = 0.3454737251382253
$ws.Range("J2").Value = 0.3454737251382253
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.856403666666667
$ws.Range("N2").Value = 8.569211
$ws.Range("O2").Value = 0.235832554697756
$ws.Range("P2").Value = 0.235832554697756
$ws.Range("Q2").Value = 4.021903006911222
$ws.Range("R2").Value = 36.19712706220099
$ws.Range("S2").Value = 0.08147395118029804
$ws.Range("T2").Value = 0.08147395118029804

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna2"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.408030333333333
$ws.Range("H3").Value = 4.224091
$ws.Range("I3").Value = 0.3454737251382253
$ws.Range("J3").Value = 0.3454737251382253
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.000300666666668
$ws.Range("N3").Value = 21.000902
$ws.Range("O3").Value = 0.5779641054021444
$ws.Range("P3").Value = 0.5779641054021444
$ws.Range("Q3").Value = 9.856635681120224
$ws.Range("R3").Value = 88.709721130082
$ws.Range("S3").Value = 0.1996714124894607
$ws.Range("T3").Value = 0.1996714124894607

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna2"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.408030333333333
$ws.Range("H4").Value = 4.224091
$ws.Range("I4").Value = 0.3454737251382253
$ws.Range("J4").Value = 0.3454737251382253
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.255294666666666
$ws.Range("N4").Value = 6.765884
$ws.Range("O4").Value = 0.1862033399000996
$ws.Range("P4").Value = 0.1862033399000996
$ws.Range("Q4").Value = 3.175523301271555
$ws.Range("R4").Value = 28.579709711444
$ws.Range("S4").Value = 0.06432836146846653
$ws.Range("T4").Value = 0.06432836146846654

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna2"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.015377
$ws.Range("H5").Value = 6.046131
$ws.Range("I5").Value = 0.494492045565236
$ws.Range("J5").Value = 0.4944920455652361
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.856403666666667
$ws.Range("N5").Value = 8.569211
$ws.Range("O5").Value = 0.235832554697756
$ws.Range("P5").Value = 0.235832554697756
$ws.Range("Q5").Value = 5.756730252515666
$ws.Range("R5").Value = 51.810572272641
$ws.Range("S5").Value = 0.1166173223833688
$ws.Range("T5").Value = 0.1166173223833688

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna2"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.015377
$ws.Range("H6").Value = 6.046131
$ws.Range("I6").Value = 0.494492045565236
$ws.Range("J6").Value = 0.4944920455652361
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.000300666666668
$ws.Range("N6").Value = 21.000902
$ws.Range("O6").Value = 0.5779641054021444
$ws.Range("P6").Value = 0.5779641054021444
$ws.Range("Q6").Value = 14.10824495668467
$ws.Range("R6").Value = 126.974204610162
$ws.Range("S6").Value = 0.2857986527435881
$ws.Range("T6").Value = 0.2857986527435881

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna2"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.015377
$ws.Range("H7").Value = 6.046131
$ws.Range("I7").Value = 0.494492045565236
$ws.Range("J7").Value = 0.4944920455652361
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.255294666666666
$ws.Range("N7").Value = 6.765884
$ws.Range("O7").Value = 0.1862033399000996
$ws.Range("P7").Value = 0.1862033399000996
$ws.Range("Q7").Value = 4.545268999422666
$ws.Range("R7").Value = 40.907420994804
$ws.Range("S7").Value = 0.09207607043827917
$ws.Range("T7").Value = 0.09207607043827917

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna2"
$ws.Range("C8").Value = "Epha4"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6522436666666667
$ws.Range("H8").Value = 1.956731
$ws.Range("I8").Value = 0.1600342292965385
$ws.Range("J8").Value = 0.1600342292965385
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.856403666666667
$ws.Range("N8").Value = 8.569211
$ws.Range("O8").Value = 0.235832554697756
$ws.Range("P8").Value = 0.235832554697756
$ws.Range("Q8").Value = 1.863071201026778
$ws.Range("R8").Value = 16.767640809241
$ws.Range("S8").Value = 0.03774128113408916
$ws.Range("T8").Value = 0.03774128113408915

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna2"
$ws.Range("C9").Value = "Epha4"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6522436666666667
$ws.Range("H9").Value = 1.956731
$ws.Range("I9").Value = 0.1600342292965385
$ws.Range("J9").Value = 0.1600342292965385
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.000300666666668
$ws.Range("N9").Value = 21.000902
$ws.Range("O9").Value = 0.5779641054021444
$ws.Range("P9").Value = 0.5779641054021444
$ws.Range("Q9").Value = 4.565901774595779
$ws.Range("R9").Value = 41.093115971362
$ws.Range("S9").Value = 0.09249404016909554
$ws.Range("T9").Value = 0.09249404016909554

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna2"
$ws.Range("C10").Value = "Epha4"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6522436666666667
$ws.Range("H10").Value = 1.956731
$ws.Range("I10").Value = 0.1600342292965385
$ws.Range("J10").Value = 0.1600342292965385
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.255294666666666
$ws.Range("N10").Value = 6.765884
$ws.Range("O10").Value = 0.1862033399000996
$ws.Range("P10").Value = 0.1862033399000996
$ws.Range("Q10").Value = 1.471001662800444
$ws.Range("R10").Value = 13.239014965204
$ws.Range("S10").Value = 0.02979890799335384
$ws.Range("T10").Value = 0.02979890799335384
